$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Prototype" row (5): written-pages formula now tracks through day 24
# instead of day 21, and the needed-pages target drops from 26 to 24.
$ws.Range("B5").Formula = "=SUM(I17:I24)"
$ws.Range("C5").Value = 24

# "Comparison" (row 6) and "Conclusion" (row 8) needed-pages targets
# each grow by one page.
$ws.Range("C6").Value = 4
$ws.Range("C8").Value = 2

# One more page got written on day 24 (I24).
$ws.Range("I24").Value = 2

# Move the active selection to K29 (cursor position as of this save).
$ws.Range("K29").Select()

$wb.Application.CalculateFullRebuild()
